$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("params")
$ws.Range("B2").Value = 0.47
$ws.Range("B3").Value = 0.53
